# Insert a new weekly price-record row into the "Apio" sheet.
#
# The source data is shifted down by one row starting at row 350 (i.e. a
# brand-new row is inserted at row 350, pushing the former rows 350-371
# down to 351-372), and the newly inserted row 350 is populated with a
# new weekly observation (date 44931, "Primera" quality, 25 units sold,
# prices 13000/13000/13000, average 2167).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 350; Excel shifts rows 350:371
# down to 351:372 and carries the formatting (incl. the date-cell style)
# down from the row above, matching the existing data's layout.
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(350, 1).Value2 = 4
$ws.Cells.Item(350, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(350, 3).Value2 = "Los Lagos"
$ws.Cells.Item(350, 4).Value2 = 44931
$ws.Cells.Item(350, 5).Value2 = 10
$ws.Cells.Item(350, 6).Value2 = 100112017
$ws.Cells.Item(350, 7).Value2 = "Apio"
$ws.Cells.Item(350, 8).Value2 = "Americana (o)"
$ws.Cells.Item(350, 9).Value2 = "Primera"
$ws.Cells.Item(350, 10).Value2 = 25
$ws.Cells.Item(350, 11).Value2 = 13000
$ws.Cells.Item(350, 12).Value2 = 13000
$ws.Cells.Item(350, 13).Value2 = 13000
$ws.Cells.Item(350, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(350, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(350, 16).Value2 = 2167
$ws.Cells.Item(350, 17).Value2 = 6
$ws.Cells.Item(350, 18).Value2 = "Hortaliza"
